# Automatische test-sync: 2025-07-27 19:39:50
# Append a new log row (#13) to the "Logs" sheet and update the
# "Dashboard" summary count for the "Productinformatie" category.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 15

$logs.Cells.Item($newRow, 1).Value  = "Kun je mij de datasheet van de VentiQ-250 sturen?"
$logs.Cells.Item($newRow, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value  = "Testmail #13: Kun je mij de datasheet van de VentiQ-250 sturen?"
$logs.Cells.Item($newRow, 4).Value  = "Productinformatie"
$logs.Cells.Item($newRow, 5).Value  = "Beste heer/mevrouw,`r`nBedankt voor uw interesse in de VentiQ-250. Helaas kan ik u op dit moment geen datasheet van de VentiQ-250 sturen, omdat ik geen e-mailbijlage kan openen.`r`nU kunt de datasheet van de VentiQ-250 vinden op onze website [www.bedrijfsnaam.nl] onder het tabblad 'Producten'. Mocht u nog verdere vragen hebben, dan help ik u graag verder.`r`nMet vriendelijke groet,`r`n[Naam]`r`nE-mailassistent"
$logs.Cells.Item($newRow, 6).Value  = "2025-07-27 19:39:08"
$logs.Cells.Item($newRow, 7).Value  = "Ja"
$logs.Cells.Item($newRow, 8).Value  = "Nee"
$logs.Cells.Item($newRow, 9).Value  = "Ja"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Setting a multi-line value auto-expands the row height; restore the
# default (unset) row height so the new row matches the others.
$logs.Rows.Item($newRow).AutoFit()

# Extend the conditional-formatting ranges so they cover the new row too
# (D2:D14 -> D2:D15, G2:G14 -> G2:G15, H2:H14 -> H2:H15, I2:I14 -> I2:I15,
# J2:J14 -> J2:J15).
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "14")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "15")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary: "Productinformatie" count 3 -> 4
$dashboard.Cells.Item(3, 2).Value = 4
